$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1581305064958239
$ws.Range("D2").Value = 0.8757964583428528

$ws.Range("C3").Value = -0.3740981916564439
$ws.Range("D3").Value = 0.7119125877667059

$ws.Range("C4").Value = 1.490169631537824
$ws.Range("D4").Value = 0.1503781433732672

$ws.Range("C5").Value = 0.4778140403806937
$ws.Range("D5").Value = 0.6374958980720455

$ws.Range("C6").Value = -0.1379889517784247
$ws.Range("D6").Value = 0.8915043528864448

$ws.Range("C7").Value = 1.712320696094954
$ws.Range("D7").Value = 0.100899240245691

$ws.Range("C8").Value = 0.5421388377919437
$ws.Range("D8").Value = 0.5931723314261372

$ws.Range("C9").Value = 2.343377191714163
$ws.Range("D9").Value = 0.02855515918675722

$ws.Range("C10").Value = 1.034220857718455
$ws.Range("D10").Value = 0.3122684366282888

$ws.Range("C11").Value = -1.282914624190289
$ws.Range("D11").Value = 0.2128736151774504
